$d = $word.ActiveDocument

# --- Change 1: merge " versione definitiva." + bookmark + " " into
#     " versione definitiva. " (drops the old _GoBack bookmark, which
#     gets re-created later near the deleted table row). ---
$d.Content.Find.Execute(" versione definitiva. ", $true, $false, $false, $false, $false, $true, 1, $false, " versione definitiva. ", 2) | Out-Null

# --- Change 2: delete the "Tipo (chiave primaria di Mazzo)" row from the
#     first table (UtenteRegistrato). ---
$tbl = $d.Tables.Item(1)
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $row = $tbl.Rows.Item($i)
    if ($row.Cells.Item(1).Range.Text -like "*Tipo (chiave primaria di Mazzo)*") {
        $row.Delete()
        break
    }
}

# --- Change 3: re-add the _GoBack bookmark to the (now empty) paragraph
#     that immediately follows the table. Use a Find-derived, non-collapsed
#     Range (matching the paragraph mark via "^p") so the bookmark lands in
#     the right spot. ---
$rng = $d.Content
$rng.Start = 4580
$rng.Find.Execute("^p", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null

# --- Change 4: Username (chiave primaria di utenteModeratore) ->
#     Username (chiave primaria di utenteRegistrato) ---
$d.Content.Find.Execute(" (chiave primaria di utenteModeratore)", $true, $false, $false, $false, $false, $true, 1, $false, " (chiave primaria di utenteRegistrato)", 2) | Out-Null

# --- Change 5: Realizza (relazione tra utenteModeratore e Storia) ->
#     Realizza (relazione tra utenteRegistrato e Storia) ---
$d.Content.Find.Execute("Realizza (relazione tra utenteModeratore e Storia)", $true, $false, $false, $false, $false, $true, 1, $false, "Realizza (relazione tra utenteRegistrato e Storia)", 2) | Out-Null

# --- Change 6: Crea (relazione tra utenteGiocatore e Personaggio) ->
#     Crea (relazione tra utenteRegistrato e Personaggio) ---
$d.Content.Find.Execute("Crea (relazione tra utenteGiocatore e Personaggio)", $true, $false, $false, $false, $false, $true, 1, $false, "Crea (relazione tra utenteRegistrato e Personaggio)", 2) | Out-Null

# --- Change 7: Ruolo (chiave primaria di utenteGiocatore) ->
#     Username (chiave primaria di utenteRegistrato) ---
$d.Content.Find.Execute("Ruolo (chiave primaria di utenteGiocatore)", $true, $false, $false, $false, $false, $true, 1, $false, "Username (chiave primaria di utenteRegistrato)", 2) | Out-Null
